# Update the app_identifier data value in the "Global" data table for the
# IOS row (row 4) from "com.hpe.iShopping" to "com.mf.iShopping".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")
$ws.Range("E4").Value = "com.mf.iShopping"
